# Update cryptocurrency price/volume data per upstream scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.299.23"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.839.01"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.93"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6240"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07372"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.71"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "1.834.47"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.944"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001054"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6610"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.19"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.239"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "29.296.59"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.72"
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.292"
$ws.Range("E22").Value = "  -3.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.39"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.414"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1334"
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07094"
$ws.Range("E28").Value = "  +8.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.484"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.024"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.016"
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.148"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.797"
$ws.Range("E34").Value = "  -3.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6936"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.585"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01824"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.233.31"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.780"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.777"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9517"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "RocketPoolETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D43").Value = "1.991.59"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000116"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.925"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.676"
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.902"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1128"
$ws.Range("E50").Value = "  -3.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3866"
$ws.Range("E51").Value = "  -2.27%  "
